# Auto-generated edit script: update cryptos price/volume data per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.309.12'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '2.948.54'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''565.87'
$ws.Range('E5').Value = '  -2.68%  '
$ws.Range('D6').Value = '''158.29'
$ws.Range('E6').Value = '  +4.12%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''0.519'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').Value = '2.943.30'
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').Value = '''6.77'
$ws.Range('E10').Value = '  -2.65%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '''0.458'
$ws.Range('E12').Value = '  +2.54%  '
$ws.Range('E13').Value = '  +2.89%  '
$ws.Range('D14').Value = '''33.96'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').Value = '65.534.47'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').Value = '3.438.77'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').Value = '''6.95'
$ws.Range('D19').Value = '2.949.34'
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').Value = '''446.28'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '''13.79'
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '''7.18'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('D24').Value = '''82.70'
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '''12.08'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').Value = '''2.17'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '''9.90'
$ws.Range('E28').Value = '  -6.33%  '
$ws.Range('D29').Value = '''7.88'
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('D30').Value = '''2.33'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '''2.56'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').Value = '0.0₃0969'
$ws.Range('E32').Value = '  -4.47%  '
$ws.Range('D33').Value = '''27.27'
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('E37').Value = '  +1.31%  '
$ws.Range('D38').Value = '''49.06'
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('E39').Value = '  -4.76%  '
$ws.Range('D40').Value = '''0.297'
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.118'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('B42').Value = 'Arweave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D42').Value = '''42.70'
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('D43').Value = '''8.43'
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('D44').Value = '''2.76'
$ws.Range('E44').Value = '  -3.77%  '
$ws.Range('D45').Value = '''383.94'
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').Value = '''0.0352'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').Value = '2.728.57'
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('D48').Value = '''130.67'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').Value = '''0.106'
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('D51').Value = '''2.14'
$ws.Range('E51').Value = '  +5.69%  '
